# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 00:26"

# --- Re-sorted country rows: swap the country names that changed rank ---
# Barein overtook Rumania
$ws.Range("A50").Value = "Barein"
$ws.Range("A51").Value = "Rumania"

# Liberia overtook Republica del Chad
$ws.Range("A143").Value = "Liberia"
$ws.Range("A144").Value = "Republica del Chad"

# Dominica overtook Fiyi
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# --- Updated case numbers ---
# columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

$updates = @{
    4   = @(2976185, 40415, 1285147, 1558492, 0, 228, 132546)
    5   = @(1603055, 24679, 978615,  559573,  0, 502, 64867)
    22  = @(117110,  3721,  47881,   65165,   0, 122, 4064)
    50  = @(29367,   510,   24649,   4621,    0, 1,   97)
    51  = @(28973,   391,   20026,   7197,    0, 19,  1750)
    70  = @(10772,   310,   5067,    5631,    0, 2,   74)
    87  = @(5740,    63,    2915,    2579,    0, 5,   246)
    111 = @(2330,    27,    1527,    684,     0, 1,   119)
    143 = @(874,     5,     377,     460,     0, 0,   37)
    144 = @(872,     1,     787,     11,      0, 0,   74)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

# Rows where only D and E (Casos activos / Recuperados) plus G/H changed
$ws.Cells.Item(32, 4).Value = 28722
$ws.Cells.Item(32, 5).Value = 28032
$ws.Cells.Item(32, 7).Value = 12
$ws.Cells.Item(32, 8).Value = 4781

# Rows where only B-E changed (F, G, H unchanged)
$ws.Cells.Item(124, 2).Value = 1542
$ws.Cells.Item(124, 3).Value = 9
$ws.Cells.Item(124, 4).Value = 1062
$ws.Cells.Item(124, 5).Value = 418

$ws.Cells.Item(133, 2).Value = 1105
$ws.Cells.Item(133, 3).Value = 13
$ws.Cells.Item(133, 4).Value = 567
$ws.Cells.Item(133, 5).Value = 535
